$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the style from row 4 (A4) to A5 so the new row matches existing formatting
$ws.Range("A4").Copy() | Out-Null
$ws.Range("A5").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Set the values for the new row 5
$ws.Range("A5").Value = 45862.57291666666
$ws.Range("B5").Value = 2025
$ws.Range("C5").Value = 30
$ws.Range("D5").Value = 18.78
$ws.Range("E5").Value = 75.75
$ws.Range("F5").Value = 91.41
$ws.Range("G5").Value = 14.3
$ws.Range("H5").Value = "ESE"
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = "14:00:12"
